$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.955270666666667
$ws.Range("N2").Value = 5.865812
$ws.Range("O2").Value = 0.4478934208563147
$ws.Range("P2").Value = 0.4478934208563147
$ws.Range("Q2").Value = 0.090356968048
$ws.Range("R2").Value = 0.8132127124320001
$ws.Range("S2").Value = 0.4478934208563147
$ws.Range("T2").Value = 0.4478934208563147

# Row 3
$ws.Range("O3").Value = 0.06160472848893509
$ws.Range("P3").Value = 0.06160472848893508
$ws.Range("S3").Value = 0.06160472848893509
$ws.Range("T3").Value = 0.06160472848893508

# Row 4
$ws.Range("M4").Value = 1.050406333333333
$ws.Range("N4").Value = 3.151219
$ws.Range("O4").Value = 0.2406163473663007
$ws.Range("P4").Value = 0.2406163473663007
$ws.Range("Q4").Value = 0.048541377476
$ws.Range("R4").Value = 0.436872397284
$ws.Range("S4").Value = 0.2406163473663007
$ws.Range("T4").Value = 0.2406163473663007

# Row 5
$ws.Range("M5").Value = 0.5333156666666666
$ws.Range("N5").Value = 1.599947
$ws.Range("O5").Value = 0.1221665022709214
$ws.Range("P5").Value = 0.1221665022709214
$ws.Range("Q5").Value = 0.024645583588
$ws.Range("R5").Value = 0.221810252292
$ws.Range("S5").Value = 0.1221665022709214
$ws.Range("T5").Value = 0.1221665022709214

# Row 6
$ws.Range("M6").Value = 0.5575549999999999
$ws.Range("N6").Value = 1.672665
$ws.Range("O6").Value = 0.1277190010175279
$ws.Range("P6").Value = 0.1277190010175279
$ws.Range("Q6").Value = 0.02576573166
$ws.Range("R6").Value = 0.23189158494
$ws.Range("S6").Value = 0.1277190010175279
$ws.Range("T6").Value = 0.1277190010175279
